$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data rows 1..542 (header-less, data starts row 1).
# We append 8 new rows (543..550) replicating the existing A/B/C pattern:
#   column A: test/scenario name, column B: status, column C: browser
# New rows reuse existing text values already present in the sheet.

$lastRow = 542

$newRows = @(
    @("Create and Delete CitizenShip From Excel", "PASSED", "chrome"),
    @("Create Country", "PASSED", "chrome"),
    @("Create Nationality", "PASSED", "chrome"),
    @("Fee Functionality", "PASSED", "chrome"),
    @("Fee Functionality", "PASSED", "chrome"),
    @("Fee Functionality", "PASSED", "chrome"),
    @("Fee Functionality", "PASSED", "chrome"),
    @("Fee Functionality", "PASSED", "chrome")
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $lastRow + 1 + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
